$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Act" path entry as row 25 (mirrors the existing "A2" row pattern at row 18)
$ws.Range("A25").Value = "Act"
$ws.Range("E25").Value = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Act"

# Update the view: change the frozen pane's top-left cell and the active selection
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("K32").Select()
